$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "Table2" currently spans A1:E66 (65 data rows + header).
# Add one new row to the table, which both extends the table range
# and the worksheet data region to A1:E67.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()
$newRow = $newListRow.Range

# Populate the new row (row 67) with the new LeetCode entry.
$newRow.Cells.Item(1, 1).Value = "62. Unique Paths"
$newRow.Cells.Item(1, 2).Value = "Medium"
$newRow.Cells.Item(1, 3).Value = "Dynamic Programming"
$newRow.Cells.Item(1, 4).Value = "A recursive solution is DFS with a cache[r][c], as res = right + down counts. For the base case, we can choose to define the finish as 1, and define out of bounds as 0. The last row and col should be filled with 1 as there is only 1 possible path. We can set each cell from the base case (end) and fill the cells with the number of paths that can reach the end. You can return the sum of all of them as the result."
$newRow.Cells.Item(1, 5).Value = "https://leetcode.com/problems/unique-paths/solutions/182143/recursive-memoization-and-dynamic-programming-solutions/ "

# Match existing formatting conventions used throughout the table:
# column B (Difficulty) carries a fill highlight per-difficulty ("Medium"
# is the orange FFC000 fill, same as the other Medium rows).
$newRow.Cells.Item(1, 2).Interior.Color = 49407

# Turn the new Link cell into a real hyperlink, like the others in the column,
# then (re)apply the built-in Hyperlink cell style so it renders the same way
# as every other link in column E.
$linkUrl = "https://leetcode.com/problems/unique-paths/solutions/182143/recursive-memoization-and-dynamic-programming-solutions/ "
$ws.Hyperlinks.Add($newRow.Cells.Item(1, 5), $linkUrl, "", "", $linkUrl) | Out-Null
$newRow.Cells.Item(1, 5).Style = "Hyperlink"

# Update the view: scrolled one row further, selection moved to D71.
$ws.Range("D71").Select() | Out-Null
